# Estado de Cuenta (NIT-9017257695) update:
#  - A new mora period (2506) is added as its own row, and the existing
#    data row is advanced to the new period (2507).
#  - Totals (Valor Mora, Cant. Periodos) are updated accordingly.
#  - The signature block at the bottom keeps its two lines, now one row
#    lower because of the inserted data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the existing worker/period data row (16),
# pushing the footer/signature rows down by one.
$ws.Rows("17:17").Insert()

# Clone the formatting (fonts, borders, number formats) of the existing
# data row into the newly inserted row.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# New row keeps the prior period (2506) data for the same worker.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "1082894538"
$ws.Range("D17").Value = "IBRAHITH JUNIELES ROSELLON"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

# The original data row now reflects the newer period (2507).
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2507"

# Update the summary figures: Valor Mora doubles (two periods in arrears)
# and Cant. Periodos goes from 1 to 2.
$ws.Range("E11").Value = 160000
$ws.Range("F13").Value = 2
